$wb = $excel.ActiveWorkbook

# "Croatia" is the last sheet and the template for each per-market tab.
$croatia = $wb.Worksheets.Item("Croatia")

# Duplicate it (equivalent to right-click > Move or Copy... > Create a copy),
# dropping the new tab immediately after "Croatia".
$croatia.Copy($null, $croatia)

# The freshly inserted sheet is now the last one in the workbook.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Fill in the market-specific values for Greece.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3205"

# "Croatia" is no longer the active tab; select all of it (as Excel does
# while copying the sheet) and make "Greece" the active tab/selection.
[void]$croatia.Cells.Select()
[void]$greece.Activate()
[void]$greece.Range("D18").Select()
